$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append (rows 552-557)
$rows = @(
    @{ A=7; B="Terminal Hortofrutícola Agro Chillán"; C="Ñuble"; D=45191; E=16; F="Fruta"; G=100104; H="Frutos de pepita"; I=100104005; J="Pera"; K="Packham's Triumph"; L="Especial"; M=60;  N=14000; O=14000; P=14000; Q="$/bandeja 18 kilos granel"; R="Región de O'Higgins"; S=778; T=18 },
    @{ A=7; B="Terminal Hortofrutícola Agro Chillán"; C="Ñuble"; D=45191; E=16; F="Fruta"; G=100104; H="Frutos de pepita"; I=100104005; J="Pera"; K="Packham's Triumph"; L="Primera";  M=60;  N=12000; O=12000; P=12000; Q="$/bandeja 18 kilos granel"; R="Región de O'Higgins"; S=667; T=18 },
    @{ A=7; B="Terminal Hortofrutícola Agro Chillán"; C="Ñuble"; D=45191; E=16; F="Fruta"; G=100104; H="Frutos de pepita"; I=100104005; J="Pera"; K="Packham's Triumph"; L="Segunda";  M=60;  N=10000; O=10000; P=10000; Q="$/bandeja 18 kilos granel"; R="Región de O'Higgins"; S=556; T=18 },
    @{ A=7; B="Terminal Hortofrutícola Agro Chillán"; C="Ñuble"; D=45191; E=16; F="Fruta"; G=100104; H="Frutos de pepita"; I=100104005; J="Pera"; K="Winter Nelis";       L="Especial"; M=60;  N=13000; O=13000; P=13000; Q="$/bandeja 18 kilos granel"; R="Región de O'Higgins"; S=722; T=18 },
    @{ A=7; B="Terminal Hortofrutícola Agro Chillán"; C="Ñuble"; D=45191; E=16; F="Fruta"; G=100104; H="Frutos de pepita"; I=100104005; J="Pera"; K="Winter Nelis";       L="Primera";  M=100; N=11000; O=11000; P=11000; Q="$/bandeja 18 kilos granel"; R="Región de O'Higgins"; S=611; T=18 },
    @{ A=7; B="Terminal Hortofrutícola Agro Chillán"; C="Ñuble"; D=45191; E=16; F="Fruta"; G=100104; H="Frutos de pepita"; I=100104005; J="Pera"; K="Winter Nelis";       L="Segunda";  M=80;  N=9000;  O=9000;  P=9000;  Q="$/bandeja 18 kilos granel"; R="Región de O'Higgins"; S=500; T=18 }
)

$startRow = 552
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data.A
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = $data.G
    $ws.Cells.Item($r, 8).Value = $data.H
    $ws.Cells.Item($r, 9).Value = $data.I
    $ws.Cells.Item($r, 10).Value = $data.J
    $ws.Cells.Item($r, 11).Value = $data.K
    $ws.Cells.Item($r, 12).Value = $data.L
    $ws.Cells.Item($r, 13).Value = $data.M
    $ws.Cells.Item($r, 14).Value = $data.N
    $ws.Cells.Item($r, 15).Value = $data.O
    $ws.Cells.Item($r, 16).Value = $data.P
    $ws.Cells.Item($r, 17).Value = $data.Q
    $ws.Cells.Item($r, 18).Value = $data.R
    $ws.Cells.Item($r, 19).Value = $data.S
    $ws.Cells.Item($r, 20).Value = $data.T
}
